# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet between "2021-Q3" and "总计"
# ---------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $sheetQ3)
$newSheet.Name = "2022-Q1"

# match the page margins used by the "总计" sheet / new sheet in the target file
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# re-resolve "2021-Q3" (index shifted because of the insertion above) so
# later Range/Copy calls definitely hit the right sheet
$sheetQ3 = $wb.Worksheets.Item("2021-Q3")

# ---- header row (bold / centered / bordered style, copied from 2021-Q3) ----
$sheetQ3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---- data rows ----
$fundData = @(
    @{Row=2; A=0; B="011550"; C="湘财创新成长一年持有期混合A";               D="2.62"; E="93.51"; F="4.57"; G="0.1197"; H=3}
    @{Row=3; A=1; B="501099"; C="平安科技创新3年封闭运作灵活配置混合型";     D="3.01"; E="91.11"; F="3.56"; G="0.1072"; H=5}
    @{Row=4; A=2; B="700004"; C="平安灵活配置混合";                         D="0.46"; E="78.72"; F="3.04"; G="0.0140"; H=10}
    @{Row=5; A=3; B="011551"; C="湘财创新成长一年持有期混合C";               D="0.28"; E="93.51"; F="4.57"; G="0.0128"; H=3}
    @{Row=6; A=4; B="010076"; C="湘财长弘灵活配置混合A";                     D="0.32"; E="91.72"; F="3.77"; G="0.0121"; H=9}
    @{Row=7; A=5; B="010307"; C="西藏东财信息产业精选混合A";                 D="0.29"; E="77.93"; F="3.43"; G="0.0099"; H=1}
    @{Row=8; A=6; B="010077"; C="湘财长弘灵活配置混合C";                     D="0.10"; E="91.72"; F="3.77"; G="0.0038"; H=9}
    @{Row=9; A=7; B="010308"; C="西藏东财信息产业精选混合C";                 D="0.09"; E="77.93"; F="3.43"; G="0.0031"; H=1}
)

foreach ($rec in $fundData) {
    $r = $rec.Row

    # column A keeps the bold/centered/bordered "index" style (copy from 2021-Q3!A2)
    $sheetQ3.Range("A2").Copy()
    $newSheet.Range("A" + $r).PasteSpecial(-4122)
    $newSheet.Range("A" + $r).Value = $rec.A

    # B/C/D/E/F/G are plain text cells (force text so things like leading
    # zeros / decimal-look-alikes are not coerced into numbers)
    $textRange = $newSheet.Range("B" + $r + ":G" + $r)
    $textRange.NumberFormat = "@"
    $newSheet.Range("B" + $r).Value = $rec.B
    $newSheet.Range("C" + $r).Value = $rec.C
    $newSheet.Range("D" + $r).Value = $rec.D
    $newSheet.Range("E" + $r).Value = $rec.E
    $newSheet.Range("F" + $r).Value = $rec.F
    $newSheet.Range("G" + $r).Value = $rec.G
    $textRange.ClearFormats()

    # H is a genuine number, default style
    $newSheet.Range("H" + $r).Value = $rec.H
}

# ---------------------------------------------------------------------
# 2. Add a new leading row to "总计" summarising the 2022-Q1 sheet
# ---------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Rows(2).Insert()

# re-resolve once more: Insert() can move things around, make sure we are
# still pointing at "总计"
$sheetTotal = $wb.Worksheets.Item("总计")

# clear the formatting copied down onto the new row's B:D cells
$sheetTotal.Range("B2:D2").ClearFormats()

# column A keeps its bold/centered/bordered "index" style
$sheetTotal.Range("A3").Copy()
$sheetTotal.Range("A2").PasteSpecial(-4122)

$sheetTotal.Range("A2").Value = 0
$sheetTotal.Range("B2").Value = "2022-Q1"
$sheetTotal.Range("C2").Value = 8
$sheetTotal.Range("D2").Value = 0.28

# renumber the index column for the rows that shifted down
$sheetTotal.Range("A3").Value = 1
$sheetTotal.Range("A4").Value = 2

# ---------------------------------------------------------------------
# 3. Restore the originally active sheet/selection ("2021-Q2")
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Select()

